$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before the old "wynik" column (E), to hold a helper
#    "Brak walidacji..." style pass/fail indicator column.
$ws.Columns("E").Insert()

# Match the width of column D (14.28515625 raw units -> 13.5 "characters"
# as reported by Excel) and hide the helper column.
$ws.Columns("E").ColumnWidth = 13.5
$ws.Columns("E").Hidden = $true

# 2. Populate the new column E (rows 3..32) with the pass/fail helper
#    formula, referencing the (now shifted) "wynik" column F.
$ws.Range("E3").Formula = '=IF(F3="pozytywny",0,1)'
For ($r = 4; $r -le 32; $r++) {
    $ws.Range("E$r").Formula = "=IF(F$r=`"pozytywny`",0,1)"
}

# A handful of rows were overridden by hand with a literal count rather
# than the computed formula result.
$hardcoded = @{14=2; 15=2; 16=2; 17=2; 18=2; 19=2; 24=3; 25=3; 26=3; 27=3; 28=3; 29=3}
foreach ($r in $hardcoded.Keys) {
    $ws.Range("E$r").Value = $hardcoded[$r]
}

# 3. Row 28's remark (now in column G after the insert) gets a new note
#    about missing validation for numeric-as-string input.
$ws.Range("G28").Value = "Brak walidacji czy wprowadzona liczba jest stringiem"

# 4. Restore the active selection to match the edited state.
$ws.Range("F34").Select() | Out-Null
